# Update cryptos list (Price / Volume(1h) columns) plus coin-row swaps
# (rows 26/27, 37/38, 40/41), matching the
# "Updated cryptos list ... with GitHub Actions" commit.
# Price-column values that look like plain numbers (e.g. "336.45", "1.000")
# are written with a leading apostrophe so Excel keeps them as literal text
# (matching the original t="inlineStr" cells) instead of silently coercing
# them to numeric doubles and losing formatting such as trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.151.39'
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').Value = '1.783.85'
$ws.Range('E3').Value = '  -1.67%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = '''336.45'
$ws.Range('E5').Value = '  -1.88%  '
$ws.Range('D6').Value = '''1.000'
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('D7').Value = '''0.3830'
$ws.Range('E7').Value = '  +0.58%  '
$ws.Range('E8').Value = '  -2.33%  '
$ws.Range('D9').Value = '''47.96'
$ws.Range('E9').Value = '  -1.58%  '
$ws.Range('E10').Value = '  -3.54%  '
$ws.Range('D11').Value = '''0.07448'
$ws.Range('E11').Value = '  -3.60%  '
$ws.Range('D12').Value = '''1.001'
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('D13').Value = '''21.68'
$ws.Range('E13').Value = '  -1.60%  '
$ws.Range('D14').Value = '''6.435'
$ws.Range('D15').Value = '1.780.19'
$ws.Range('E15').Value = '  -1.61%  '
$ws.Range('D16').Value = '''7.091'
$ws.Range('E16').Value = '  -2.01%  '
$ws.Range('E17').Value = '  -2.34%  '
$ws.Range('D18').Value = '''0.06640'
$ws.Range('E18').Value = '  -1.28%  '
$ws.Range('E19').Value = '  -2.89%  '
$ws.Range('D20').Value = '''1.001'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('D21').Value = '''6.519'
$ws.Range('E21').Value = '  -0.56%  '
$ws.Range('D22').Value = '''17.37'
$ws.Range('E22').Value = '  -1.27%  '
$ws.Range('D23').Value = '27.151.76'
$ws.Range('E23').Value = '  -0.83%  '
$ws.Range('E24').Value = '  -7.86%  '
$ws.Range('D25').Value = '''2.384'
$ws.Range('E25').Value = '  -3.67%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').Value = '''2.497'
$ws.Range('E26').Value = '  -6.26%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''21.12'
$ws.Range('E27').Value = '  -3.91%  '
$ws.Range('D28').Value = '''1.443'
$ws.Range('E28').Value = '  -2.03%  '
$ws.Range('D29').Value = '''155.11'
$ws.Range('E29').Value = '  +0.68%  '
$ws.Range('D30').Value = '1.981.10'
$ws.Range('E30').Value = '  -1.40%  '
$ws.Range('D31').Value = '''134.16'
$ws.Range('E31').Value = '  -1.26%  '
$ws.Range('D32').Value = '''3.981'
$ws.Range('E32').Value = '  -1.17%  '
$ws.Range('D33').Value = '''6.040'
$ws.Range('E33').Value = '  -4.33%  '
$ws.Range('D34').Value = '''0.08652'
$ws.Range('E34').Value = '  -1.35%  '
$ws.Range('D35').Value = '''12.99'
$ws.Range('E35').Value = '  -6.38%  '
$ws.Range('D36').Value = '''1.626'
$ws.Range('E36').Value = '  -5.03%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').Value = '''5.387'
$ws.Range('E37').Value = '  -3.91%  '
$ws.Range('B38').Value = 'TheSandbox'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D38').Value = '''0.6829'
$ws.Range('E38').Value = '  -2.07%  '
$ws.Range('D39').Value = '''0.06296'
$ws.Range('E39').Value = '  -2.57%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '''0.02326'
$ws.Range('E40').Value = '  -3.75%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Value = '''0.2176'
$ws.Range('E41').Value = '  -4.23%  '
$ws.Range('E42').Value = '  -4.43%  '
$ws.Range('D43').Value = '''8.349'
$ws.Range('E43').Value = '  -6.61%  '
$ws.Range('D44').Value = '''14.21'
$ws.Range('E44').Value = '  -3.59%  '
$ws.Range('D45').Value = '''1.000'
$ws.Range('E45').Value = '  -0.13%  '
$ws.Range('D46').Value = '''0.6416'
$ws.Range('E46').Value = '  -1.52%  '
$ws.Range('D47').Value = '''3.851'
$ws.Range('E47').Value = '  -4.72%  '
$ws.Range('D48').Value = '''2.131'
$ws.Range('E48').Value = '  -2.13%  '
$ws.Range('D49').Value = '''131.29'
$ws.Range('E49').Value = '  -1.03%  '
$ws.Range('E50').Value = '  -3.17%  '
$ws.Range('D51').Value = '''78.59'
$ws.Range('E51').Value = '  -2.23%  '
